$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Row 93 already contains the note "electric death animation isn't drawing
# now for some reason" in column D. Mark it resolved ("DONE" / "Dave" /
# date) just like the other finished wish-list rows by filling in columns
# A-C (formatting copied from the row directly above, which already uses
# the correct "DONE" style).
# ---------------------------------------------------------------------------
$null = $ws.Range("A92:D92").Copy()
$null = $ws.Range("A93:D93").PasteSpecial(-4122)
$ws.Range("A93").Value = "DONE"
$ws.Range("B93").Value = "Dave"
$ws.Range("C93").Value = 39964

# ---------------------------------------------------------------------------
# Add a new finished wish-list item above the "Nifty Shit" section: a new
# row is inserted at row 160 (pushing everything below down by one) and
# filled in the same way.
# ---------------------------------------------------------------------------
$null = $ws.Rows("160:160").Insert()
$null = $ws.Range("A159:D159").Copy()
$null = $ws.Range("A160:D160").PasteSpecial(-4122)
$ws.Range("A160").Value = "DONE"
$ws.Range("B160").Value = "Dave"
$ws.Range("C160").Value = 39964
$ws.Range("D160").Value = "level select state images are glitching out"

# ---------------------------------------------------------------------------
# Restore the view: scrolled roughly to row 64 with C94 selected.
# ---------------------------------------------------------------------------
$null = $ws.Range("A64").Select()
$null = $ws.Range("C94").Select()
